$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 336, pushing the existing rows 336-346 down to 337-347.
$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with the new weekly record.
$ws.Cells.Item(336, 1).Value = 2
$ws.Cells.Item(336, 2).Value = "Comercializadora del Agro de Limar" + [char]237
$ws.Cells.Item(336, 3).Value = "Coquimbo"
$ws.Cells.Item(336, 4).Value = 44448
$ws.Cells.Item(336, 5).Value = 4
$ws.Cells.Item(336, 6).Value = 100112020
$ws.Cells.Item(336, 7).Value = "Tomate"
$ws.Cells.Item(336, 8).Value = "Larga vida"
$ws.Cells.Item(336, 9).Value = "Tercera"
$ws.Cells.Item(336, 10).Value = 1200
$ws.Cells.Item(336, 11).Value = 12000
$ws.Cells.Item(336, 12).Value = 13000
$ws.Cells.Item(336, 13).Value = 12500
$ws.Cells.Item(336, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(336, 15).Value = "Provincia de Limar" + [char]237
$ws.Cells.Item(336, 16).Value = 694
$ws.Cells.Item(336, 17).Value = 18
$ws.Cells.Item(336, 18).Value = "Hortaliza"
